$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set Week 1 (column D) remaining amounts to 0 for rows 5-7
$ws.Range("D5").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("D7").Value = 0

# Update the active cell selection to D7
$ws.Range("D7").Select()
